$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 685.86365
$ws.Range("J17").Value = 685.86365
$ws.Range("L17").Value = 2057.59095
$ws.Range("N17").Value = -2393.59095
$ws.Range("H63").Value = 163928.42
$ws.Range("I63").Value = 50000
$ws.Range("J63").Value = 182916.5
$ws.Range("K63").Value = 50000
$ws.Range("L63").Value = 182916.5
$ws.Range("M63").Value = -49376
$ws.Range("N63").Value = -184164.5
$ws.Range("H66").Value = 163928.42
$ws.Range("I66").Value = 50000
$ws.Range("J66").Value = 182916.5
$ws.Range("K66").Value = 150000
$ws.Range("L66").Value = 548749.5
$ws.Range("M66").Value = -146880
$ws.Range("N66").Value = -554989.5
$ws.Range("H68").Value = 199500
$ws.Range("J68").Value = 199500
$ws.Range("L68").Value = 199500
$ws.Range("N68").Value = -200998
$ws.Range("H69").Value = 16499.318
$ws.Range("I69").Value = 8999.333000000001
$ws.Range("K69").Value = 26997.999
$ws.Range("M69").Value = -26123.999
$ws.Range("H71").Value = 199500
$ws.Range("J71").Value = 199500
$ws.Range("L71").Value = 598500
$ws.Range("N71").Value = -605988
$ws.Range("H72").Value = 16499.318
$ws.Range("I72").Value = 8999.333000000001
$ws.Range("K72").Value = 80993.997
$ws.Range("M72").Value = -76625.997
$ws.Range("H74").Value = 19550
$ws.Range("I74").Value = 19550
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 19550
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -18614
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 19550
$ws.Range("I77").Value = 19550
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 97750
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -93070
$ws.Range("N77").ClearContents()
$ws.Range("H80").Value = 1075.6428
$ws.Range("I80").Value = 1115.5714
$ws.Range("J80").Value = 1035.7142
$ws.Range("K80").Value = 3346.7142
$ws.Range("L80").Value = 3107.1426
$ws.Range("M80").Value = -2348.7142
$ws.Range("N80").Value = -5103.142599999999
$ws.Range("H83").Value = 1075.6428
$ws.Range("I83").Value = 1115.5714
$ws.Range("J83").Value = 1035.7142
$ws.Range("K83").Value = 10040.1426
$ws.Range("L83").Value = 9321.427799999999
$ws.Range("M83").Value = -5048.142600000001
$ws.Range("N83").Value = -19305.4278
$ws.Range("H97").Value = 3099.875
$ws.Range("J97").Value = 3099.875
$ws.Range("L97").Value = 9299.625
$ws.Range("N97").Value = -10291.625
$ws.Range("H100").Value = 1246.1875
$ws.Range("I100").Value = 840.6923
$ws.Range("J100").Value = 3003.3333
$ws.Range("K100").Value = 840.6923
$ws.Range("L100").Value = 3003.3333
$ws.Range("M100").Value = -299.6923
$ws.Range("N100").Value = -4085.3333
$ws.Range("H131").Value = 1977.4
$ws.Range("I131").Value = 1977.4
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 5932.200000000001
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -892.2000000000007
$ws.Range("N131").ClearContents()
$ws.Range("H138").Value = 5515.222
$ws.Range("I138").Value = 4215.923
$ws.Range("J138").Value = 5927.1953
$ws.Range("K138").Value = 12647.769
$ws.Range("L138").Value = 17781.5859
$ws.Range("M138").Value = -7507.769
$ws.Range("N138").Value = -28061.5859

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4826.82
$ws.Range("I32").Value = 3651.5593
$ws.Range("J32").Value = 39497
$ws.Range("K32").Value = 3651.5593
$ws.Range("L32").Value = 39497
$ws.Range("M32").Value = -3364.5593
$ws.Range("N32").Value = -40071
$ws.Range("H74").Value = 4716.037
$ws.Range("I74").Value = 2415.238
$ws.Range("K74").Value = 2415.238
$ws.Range("M74").Value = -1541.238
$ws.Range("H77").Value = 4716.037
$ws.Range("I77").Value = 2415.238
$ws.Range("K77").Value = 12076.19
$ws.Range("M77").Value = -7708.189999999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2618.125
$ws.Range("I94").Value = 2050.75
$ws.Range("K94").Value = 2050.75
$ws.Range("M94").Value = -1599.75
$ws.Range("H99").Value = 2199.8572
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 3133
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 3133
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -6129

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 114880.97
$ws.Range("I7").Value = 148730.7
$ws.Range("K7").Value = 148730.7
$ws.Range("M7").Value = -148617.7
$ws.Range("H31").Value = 6419.46
$ws.Range("I31").Value = 3860.6667
$ws.Range("J31").Value = 12999.214
$ws.Range("K31").Value = 3860.6667
$ws.Range("L31").Value = 12999.214
$ws.Range("M31").Value = -3565.6667
$ws.Range("N31").Value = -13589.214
$ws.Range("H34").Value = 6419.46
$ws.Range("I34").Value = 3860.6667
$ws.Range("J34").Value = 12999.214
$ws.Range("K34").Value = 3860.6667
$ws.Range("L34").Value = 12999.214
$ws.Range("M34").Value = -3658.6667
$ws.Range("N34").Value = -13403.214
$ws.Range("H105").Value = 125002140
$ws.Range("J105").Value = 500001340
$ws.Range("L105").Value = 500001340
$ws.Range("N105").Value = -500004834
$ws.Range("H122").Value = 2075
$ws.Range("I122").Value = 1900
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 5700
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -3250
$ws.Range("N122").Value = -13750
$ws.Range("H132").Value = 4000.1
$ws.Range("I132").Value = 3365.8147
$ws.Range("J132").Value = 5317.4614
$ws.Range("K132").Value = 10097.4441
$ws.Range("L132").Value = 15952.3842
$ws.Range("M132").Value = -7567.444100000001
$ws.Range("N132").Value = -21012.3842
$ws.Range("H134").Value = 3765.1
$ws.Range("I134").Value = 2905.6667
$ws.Range("K134").Value = 8717.000100000001
$ws.Range("M134").Value = -6182.000100000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 2996.6667
$ws.Range("I97").Value = 1995
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 5985
$ws.Range("L97").Value = 15000
$ws.Range("M97").Value = -5489
$ws.Range("N97").Value = -15992
$ws.Range("H113").Value = 991.52
$ws.Range("J113").Value = 1092.6666
$ws.Range("L113").Value = 3277.9998
$ws.Range("N113").Value = -7617.9998
$ws.Range("H140").Value = 1453.4667
$ws.Range("J140").Value = 1734.0769
$ws.Range("L140").Value = 5202.2307
$ws.Range("N140").Value = -15562.2307
$ws.Range("H141").Value = 4157.5
$ws.Range("I141").Value = 1894.2858
$ws.Range("J141").Value = 20000
$ws.Range("K141").Value = 5682.857400000001
$ws.Range("L141").Value = 60000
$ws.Range("M141").Value = -502.8574000000008
$ws.Range("N141").Value = -70360

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 111555
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H80").Value = 11726.272
$ws.Range("J80").Value = 14960
$ws.Range("L80").Value = 14960
$ws.Range("N80").Value = -16956
$ws.Range("H83").Value = 11726.272
$ws.Range("J83").Value = 14960
$ws.Range("L83").Value = 74800
$ws.Range("N83").Value = -84784
$ws.Range("H97").Value = 400.66666
$ws.Range("I97").Value = 400.66666
$ws.Range("K97").Value = 400.66666
$ws.Range("M97").Value = 95.33334000000002
$ws.Range("H130").Value = 74963.71000000001
$ws.Range("J130").Value = 77516.84
$ws.Range("L130").Value = 77516.84
$ws.Range("N130").Value = -87556.84

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5856.304
$ws.Range("I7").Value = 6141.8945
$ws.Range("K7").Value = 6141.8945
$ws.Range("M7").Value = -6029.8945
$ws.Range("H43").Value = 109966.664
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 109966.664
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 109966.664
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -110352.664
$ws.Range("H100").Value = 5469.6665
$ws.Range("I100").Value = 5993.6
$ws.Range("J100").Value = 2850
$ws.Range("K100").Value = 5993.6
$ws.Range("L100").Value = 2850
$ws.Range("M100").Value = -5452.6
$ws.Range("N100").Value = -3932
$ws.Range("H126").Value = 5856.304
$ws.Range("I126").Value = 6141.8945
$ws.Range("K126").Value = 18425.6835
$ws.Range("M126").Value = -15955.6835
$ws.Range("H132").Value = 4693.591
$ws.Range("J132").Value = 6450
$ws.Range("L132").Value = 19350
$ws.Range("N132").Value = -24410
